$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 150
$ws.Range("F3").Value = 969
$ws.Range("F4").Value = 605
$ws.Range("F5").Value = 2990
$ws.Range("F6").Value = 801
$ws.Range("F7").Value = 595
$ws.Range("F8").Value = 603
$ws.Range("F9").Value = 439
$ws.Range("F12").Value = 546
$ws.Range("F14").Value = 2163
$ws.Range("F16").Value = 748
$ws.Range("F17").Value = 23
$ws.Range("F19").Value = 2681
$ws.Range("F23").Value = 533
$ws.Range("F25").Value = 645
$ws.Range("F26").Value = 15
$ws.Range("F27").Value = 24
$ws.Range("F29").Value = 11
$ws.Range("F33").Value = 124
$ws.Range("F34").Value = 904
$ws.Range("F35").Value = 4692
$ws.Range("F36").Value = 256
$ws.Range("F37").Value = 34
$ws.Range("F38").Value = 7

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 11
$ws.Range("F8").Value = 338
$ws.Range("F27").Value = 178
$ws.Range("F31").Value = 26
$ws.Range("F35").Value = 26
$ws.Range("F37").Value = 555
$ws.Range("F38").Value = 18

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 575
$ws.Range("F6").Value = 262
$ws.Range("F7").Value = 264

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 575
$ws.Range("F5").Value = 150
$ws.Range("F6").Value = 262
$ws.Range("F7").Value = 969
$ws.Range("F8").Value = 605
$ws.Range("F9").Value = 2990
$ws.Range("F10").Value = 801
$ws.Range("F11").Value = 595
$ws.Range("F12").Value = 603
$ws.Range("F13").Value = 439
$ws.Range("F16").Value = 546
$ws.Range("F17").Value = 338
$ws.Range("F21").Value = 2163
$ws.Range("F23").Value = 748
$ws.Range("F27").Value = 2681
$ws.Range("F32").Value = 533
$ws.Range("F35").Value = 264
$ws.Range("F37").Value = 645
$ws.Range("F38").Value = 645
$ws.Range("F45").Value = 904
$ws.Range("F47").Value = 4692
$ws.Range("F48").Value = 256
$ws.Range("F50").Value = 555
